$wb = $excel.ActiveWorkbook

# The small lookup sheet is named "Sheet1" (it holds JM101 / JJ101 reference codes).
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new row above the current first row, pushing the existing
# two rows (JM101, JJ101) down by one, and put "Faulty" into the new A1.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Faulty"

# Match the diff's recorded selection: the active cell ends up back on A1.
$ws.Range("A1").Select() | Out-Null
